# Test case names updated as per documentation (smoke sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("smoke")

$ws.Range("A17").Value = "AV_2268_Validate_Accuracy_of_COBT_For_DIALCelebi_User"
$ws.Range("A18").Value = "AV_2268_Validate_Accuracy_of_COBT_For_GMR_HYD_AISATS_User"
$ws.Range("A19").Value = "AV_2268_Validate_Accuracy_of_COBT_For_GMR_HYD_SG_User"

$ws.Range("A20").Value = "AV_2293_Identify_coverage_of_Flight_Sensor_and_Validate_timestamps_of_Arrival_Aircrafts"
$ws.Rows.Item(20).AutoFit()

$ws.Range("A21").Value = "AV_2294_Identify_coverage_of_Flight_Sensor_and_Validate_timestamps_of_Departure_Aircrafts"
$ws.Rows.Item(21).AutoFit()

$ws.Range("A22").Value = "AV_2307_Validate_LANDING_ONBLOCK_OFFBLOCK_AIRBORNE_timestamps_of_Arrival_and_Departure_aircrafts_Any_Data_source"
$ws.Range("B22").Value = "N"

$ws.Range("A23").Value = "AV_2304_Identify_the_coverage_of_Boarding_activities_and_validate_timestamps"
$ws.Range("B23").Value = "N"
$ws.Range("B22").Copy()
$ws.Range("B23").PasteSpecial(-4122)

$ws.Activate()
$ws.Range("A24").Select()
